$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format price cells that would otherwise be auto-parsed as numbers,
# so their exact original text (including trailing zeros) is preserved.
$textRows = @(5,9,19,20,23,25,27,30,36,39,42,43,46,47,49,50)
foreach ($r in $textRows) {
    $ws.Cells.Item($r, 4).NumberFormat = "@"
}

# Apply updated values
$ws.Range('D2').Value = '26.176.84'
$ws.Range('E2').Value = '  -0.79%  '
$ws.Range('D3').Value = '1.585.18'
$ws.Range('E3').Value = '  -0.53%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').Value = '211.71'
$ws.Range('E5').Value = '  +0.76%  '
$ws.Range('E6').Value = '  -0.54%  '
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('E8').Value = '  -0.75%  '
$ws.Range('D9').Value = '0.0602'
$ws.Range('E9').Value = '  -1.54%  '
$ws.Range('E10').Value = '  -2.35%  '
$ws.Range('E11').Value = '  +0.22%  '
$ws.Range('D12').Value = '1.808.76'
$ws.Range('E12').Value = '  -0.44%  '
$ws.Range('D13').Value = '1.582.64'
$ws.Range('E13').Value = '  -0.36%  '
$ws.Range('E14').Value = '  -1.96%  '
$ws.Range('E15').Value = '  -0.45%  '
$ws.Range('E16').Value = '  -1.24%  '
$ws.Range('D17').Value = '26.179.13'
$ws.Range('E17').Value = '  -0.76%  '
$ws.Range('D18').Value = '0.0₃0722'
$ws.Range('E18').Value = '  -0.99%  '
$ws.Range('D19').Value = '213.94'
$ws.Range('E19').Value = '  +1.05%  '
$ws.Range('D20').Value = '7.26'
$ws.Range('E20').Value = '  -2.95%  '
$ws.Range('E21').Value = '  -0.07%  '
$ws.Range('E22').Value = '  -0.82%  '
$ws.Range('D23').Value = '8.94'
$ws.Range('E23').Value = '  +0.11%  '
$ws.Range('E24').Value = '  -2.62%  '
$ws.Range('D25').Value = '144.13'
$ws.Range('E25').Value = '  -0.86%  '
$ws.Range('E26').Value = '  +0.01%  '
$ws.Range('D27').Value = '6.96'
$ws.Range('E27').Value = '  -1.41%  '
$ws.Range('E28').Value = '  -1.33%  '
$ws.Range('E29').Value = '  -1.43%  '
$ws.Range('D30').Value = '0.0494'
$ws.Range('E30').Value = '  -2.17%  '
$ws.Range('E31').Value = '  +0.28%  '
$ws.Range('E32').Value = '  -1.31%  '
$ws.Range('D33').Value = '1.408.26'
$ws.Range('E33').Value = '  +7.69%  '
$ws.Range('E34').Value = '  -2.01%  '
$ws.Range('E35').Value = '  -0.62%  '
$ws.Range('D36').Value = '0.590'
$ws.Range('E36').Value = '  -3.84%  '
$ws.Range('E37').Value = '  -1.52%  '
$ws.Range('E38').Value = '  -1.55%  '
$ws.Range('D39').Value = '0.819'
$ws.Range('E39').Value = '  +0.67%  '
$ws.Range('E40').Value = '  +4.04%  '
$ws.Range('E41').Value = '  +0.00%  '
$ws.Range('D42').Value = '0.939'
$ws.Range('E42').Value = '  -15.74%  '
$ws.Range('D43').Value = '0.764'
$ws.Range('E43').Value = '  +0.01%  '
$ws.Range('E44').Value = '  -0.24%  '
$ws.Range('D45').Value = '1.720.29'
$ws.Range('E45').Value = '  -0.49%  '
$ws.Range('D46').Value = '60.95'
$ws.Range('E46').Value = '  -2.90%  '
$ws.Range('D47').Value = '85.54'
$ws.Range('E47').Value = '  -2.98%  '
$ws.Range('E48').Value = '  -1.25%  '
$ws.Range('D49').Value = '0.0500'
$ws.Range('E49').Value = '  -1.03%  '
$ws.Range('D50').Value = '0.0968'
$ws.Range('E50').Value = '  -1.79%  '
$ws.Range('E51').Value = '  +0.01%  '
